$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: locate the start offset of $needle, searching forward from
# character offset $from in the document.
# ------------------------------------------------------------------
function Find-Pos([int]$from, [string]$needle) {
    $r = $d.Range($from, $d.Content.End)
    $ok = $r.Find.Execute($needle, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find '$needle' from offset $from"
    }
    return $r.Start
}

# Helper: split the paragraph's run at character offset $pos without
# touching its text -- drop a throw-away bookmark there and remove it
# again. Word always breaks a run around a bookmark, so this leaves
# two runs behind once the bookmark itself is gone.
function Split-At([int]$pos) {
    $rng = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplit", $rng) | Out-Null
    $d.Bookmarks("TmpSplit").Delete()
}

# ------------------------------------------------------------------
# 1) Fix the line: remove "dentro suyo " from the last sentence.
#    "Dentro de cada archivo se pueden encontrar dentro suyo más
#    elementos que pertenecen a esa rama."
#    ->
#    "Dentro de cada archivo se pueden encontrar más elementos que
#    pertenecen a esa rama."
# ------------------------------------------------------------------
$d.Content.Find.Execute("dentro suyo ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Re-type over the seam that used to separate "...el proyecto; d"
#    from "entro de esta también se encuentra la" so the two runs
#    collapse back into a single run.
# ------------------------------------------------------------------
$seam = "entro de esta también se encuentra la"
$d.Content.Find.Execute($seam, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $seam, 2) | Out-Null

# ------------------------------------------------------------------
# 3) The previous step's side effect is that every run from there to
#    the end of the paragraph (they all share the same formatting)
#    gets coalesced into one big run. Put back the original run
#    boundaries for everything after the merged "...encuentra la"
#    run, up to (but not including) the final sentence, which is
#    handled separately in step 4.
# ------------------------------------------------------------------
$p1 = Find-Pos 0 "s estructura y regla"     # before "s"
Split-At $p1

$p2 = $p1 + 1                                # after "s" / before " estructura y regla"
Split-At $p2

$p3 = Find-Pos $p2 "s de nombrado permi"     # before second "s"
Split-At $p3

$p4 = $p3 + 1                                # after "s" / before " de nombrado permi"
Split-At $p4

$p5 = Find-Pos $p4 "tiendo"                  # before "tiendo"
Split-At $p5

$p6 = $p5 + ("tiendo").Length                # after "tiendo" / before " un "
Split-At $p6

$p7 = Find-Pos $p6 "simple"                  # before "simple"
Split-At $p7

$p8 = $p7 + ("simple").Length                # after "simple" / before " e intuitivo..."
Split-At $p8

$p9 = Find-Pos $p8 "Dentro de cada archivo"  # before "Dentro de cada archivo"
Split-At $p9

# ------------------------------------------------------------------
# 4) Drop the real "_GoBack" bookmark where the edit happened, right
#    after "... se pueden encontrar ", splitting that final run in
#    two, matching what Word leaves behind after the last edit point.
# ------------------------------------------------------------------
$goBackPos = Find-Pos $p9 "encontrar "
$goBackPos = $goBackPos + ("encontrar ").Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
